$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71-192 down to 72-193.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly data point.
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44540
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112012
$ws.Cells.Item(71, 7).Value = "Espinaca"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 3200
$ws.Cells.Item(71, 11).Value = 400
$ws.Cells.Item(71, 12).Value = 500
$ws.Cells.Item(71, 13).Value = 450
$ws.Cells.Item(71, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(71, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(71, 16).Value = 900
$ws.Cells.Item(71, 17).Value = 0.5
$ws.Cells.Item(71, 18).Value = "Hortaliza"
